$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "22.471.71"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.571.07"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "290.67"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3708"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.59%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "50.03"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.90%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3372"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.66%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.146"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07530"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("E12").Value = "  -0.09%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.15"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("E14").Value = "  +0.80%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.961"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.571.49"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001119"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "90.44"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06766"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  -0.10%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.344"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "16.40"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("E23").Value = "  +2.37%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "22.482.06"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.368"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.44%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.613"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.58%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.01"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "149.07"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.075"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "125.09"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.749.75"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.067"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.94%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.186"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.015"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.798"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.69%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.08345"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.36%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02479"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.362"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.18%  "

$ws.Range("E39").Value = "  +0.60%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.06538"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("E41").Value = "  +0.34%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.28"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6220"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "14.11"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("E45").Value = "  -0.07%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.804"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("E47").Value = "  -1.29%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "129.13"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.61%  "

$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("E51").Value = "  -0.28%  "
